$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $escaped = $value -replace '"', '""'
    $c.Formula = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextCell 'D2' '67.730.51'
$ws.Range('E2').Value = '  -1.00%  '
Set-TextCell 'D3' '3.790.89'
$ws.Range('E3').Value = '  +1.20%  '
Set-TextCell 'D4' '1.00'
$ws.Range('E4').Value = '  +0.05%  '
Set-TextCell 'D5' '595.92'
$ws.Range('E5').Value = '  +0.56%  '
Set-TextCell 'D6' '166.94'
$ws.Range('E6').Value = '  +0.65%  '
Set-TextCell 'D7' '3.791.45'
$ws.Range('E7').Value = '  +1.30%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -0.02%  '
Set-TextCell 'D10' '0.159'
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('E11').Value = '  -1.77%  '
Set-TextCell 'D12' '0.450'
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('E13').Value = '  -2.16%  '
Set-TextCell 'D14' '36.00'
$ws.Range('E14').Value = '  +0.21%  '
Set-TextCell 'D15' '4.426.05'
$ws.Range('E15').Value = '  +1.21%  '
Set-TextCell 'D16' '3.776.41'
$ws.Range('E16').Value = '  +0.93%  '
Set-TextCell 'D17' '18.47'
$ws.Range('E17').Value = '  +3.20%  '
Set-TextCell 'D18' '67.739.55'
Set-TextCell 'D19' '7.03'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('E20').Value = '  -0.18%  '
Set-TextCell 'D21' '10.03'
$ws.Range('E21').Value = '  -5.92%  '
Set-TextCell 'D22' '460.24'
$ws.Range('E22').Value = '  -0.64%  '
Set-TextCell 'D23' '0.697'
$ws.Range('E23').Value = '  +0.43%  '
Set-TextCell 'D24' '0.0000154'
$ws.Range('E24').Value = '  +5.65%  '
Set-TextCell 'D25' '83.38'
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('E27').Value = '  -1.58%  '
Set-TextCell 'D28' '10.03'
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('E29').Value = '  +0.21%  '
Set-TextCell 'D30' '3.934.98'
$ws.Range('E30').Value = '  +1.14%  '
$ws.Range('E31').Value = '  +0.70%  '
Set-TextCell 'D32' '2.23'
$ws.Range('E32').Value = '  +4.37%  '
Set-TextCell 'D33' '7.21'
$ws.Range('E33').Value = '  -0.89%  '
Set-TextCell 'D34' '29.67'
$ws.Range('E34').Value = '  -0.59%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell 'D35' '0.999'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D36' '9.10'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('E39').Value = '  +0.82%  '
Set-TextCell 'D40' '0.997'
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  -0.03%  '
Set-TextCell 'D43' '46.16'
$ws.Range('E43').Value = '  +6.34%  '
$ws.Range('E44').Value = '  -0.01%  '
Set-TextCell 'D45' '48.14'
$ws.Range('E45').Value = '  +3.43%  '
$ws.Range('E46').Value = '  -0.38%  '
Set-TextCell 'D47' '149.34'
$ws.Range('E47').Value = '  +3.34%  '
$ws.Range('E48').Value = '  -1.45%  '
Set-TextCell 'D49' '394.44'
$ws.Range('E49').Value = '  +1.66%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D50' '26.72'
$ws.Range('E50').Value = '  +4.38%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell 'D51' '1.82'
$ws.Range('E51').Value = '  -4.55%  '
